$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "start", "end", "width" and "strand" columns (B:E) are removed from
# the header row and the three gene-data rows beneath it. The remaining
# columns - logFC, p-value, FDR, add.variable - currently sitting in F:I
# for those same four rows, shift left into B:E. Rows 5-16 only ever held
# empty pre-formatted F/G cells (no B:E data), so they are left as-is.
$ws.Range("F1:I4").Cut($ws.Range("B1"))

# Cutting a rectangular block also stamps blank cells for the positions
# that had no source data (D2:E4) and leaves empty-but-styled husks behind
# at the old location (F1:I4); clear both away so the saved XML has no
# stray empty <c> nodes, matching a plain column delete.
$ws.Range("D2:E4").Clear()
$ws.Range("F1:I4").Clear()

# The selection Excel stored in the file when it was last saved.
$ws.Range("C9").Select()
